$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week label bumped from "Tuần 1" to "Tuần 2"
$ws.Range("B6").Value = "Tuần 2"

# All three task rows now share the same updated task description
$ws.Range("B11").Value = "Tham gia trao đổi và hoàn thành nội dung Project Charter"
$ws.Range("B12").Value = "Tham gia trao đổi và hoàn thành nội dung Project Charter"
$ws.Range("B13").Value = "Tham gia trao đổi và hoàn thành nội dung Project Charter"

# Match the saved selection state from the workbook
$ws.Range("E19").Select() | Out-Null
